$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.126245260238647
$ws.Range("B1").Value = 4.244960308074951
$ws.Range("C1").Value = 4.248246669769287
$ws.Range("D1").Value = 1.917420148849487
$ws.Range("E1").Value = 1.246058821678162
